$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-21 with the reordered language/value data
$ws.Range("A2").Value = "English"
$ws.Range("B2").Value = 26.95727738581348
$ws.Range("A3").Value = "Chinese"
$ws.Range("B3").Value = 8.966387231973854
$ws.Range("A4").Value = "Spanish"
$ws.Range("B4").Value = 7.842536128168089
$ws.Range("A5").Value = "Japanese"
$ws.Range("B5").Value = 6.747118904652358
$ws.Range("A6").Value = "German"
$ws.Range("B6").Value = 5.780667612307163
$ws.Range("A7").Value = "Arabic"
$ws.Range("B7").Value = 4.710817813268215
$ws.Range("A8").Value = "Portuguese"
$ws.Range("B8").Value = 3.649009595082675
$ws.Range("A9").Value = "French"
$ws.Range("B9").Value = 3.502899762290265
$ws.Range("A10").Value = "Italian"
$ws.Range("B10").Value = 3.277338429977486
$ws.Range("A11").Value = "Russian"
$ws.Range("B11").Value = 3.241974760402772
$ws.Range("A12").Value = "Malay-Indonesian"
$ws.Range("B12").Value = 2.568559960101053
$ws.Range("A13").Value = "Dutch"
$ws.Range("B13").Value = 1.642638008360786
$ws.Range("A14").Value = "Korean"
$ws.Range("B14").Value = 1.599688441211686
$ws.Range("A15").Value = "Persian"
$ws.Range("B15").Value = 1.410581771220485
$ws.Range("A16").Value = "Turkish"
$ws.Range("B16").Value = 1.276164657919602
$ws.Range("A17").Value = "Thai"
$ws.Range("B17").Value = 0.9184230310573114
$ws.Range("A18").Value = "Polish"
$ws.Range("B18").Value = 0.8824618680545053
$ws.Range("A19").Value = "Urdu"
$ws.Range("B19").Value = 0.8006071592078532
$ws.Range("A20").Value = "Swedish"
$ws.Range("B20").Value = 0.5121123735386724
$ws.Range("A21").Value = "Bengali"
$ws.Range("B21").Value = 0.4388589229337661

# Remove the now-unused trailing rows (previously Uzbek / Vietnamese)
$ws.Rows(23).Delete()
$ws.Rows(22).Delete()
